$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = 123
$ws.Range("C5").Value = 123
$ws.Range("A1").Formula = "=SUM(A5,C5)"

$ws.Range("H10").Select()
